# Daily "refresh" pass over the tracker sheet.
#
# Columns: D = total days (总天), E = days remaining (剩余), F = start date
# (开始时间, yyyyMMdd). Each row's implicit deadline is F + D days; the
# sheet is regenerated once a day by recomputing E as (deadline - today).
# That means a normal row just loses one day of "remaining" (E -> E-1).
# Once a row's remaining count would hit zero, the contract/order is
# treated as renewed starting "today": F is bumped to today's date and E
# is reset back to the full D.
#
# This run's reference date is 2026-02-17 (one day after the previous
# pass, which used 2026-02-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20260217

# Data rows start at 2 (row 1 is the header) and run to the sheet's last
# used row.
$firstRow = 2
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt $firstRow) { $lastRow = 99 }

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $totalDays = $ws.Cells.Item($r, 4).Value2   # D: 总天
    $remaining = $ws.Cells.Item($r, 5).Value2   # E: 剩余
    $startDate = $ws.Cells.Item($r, 6).Value2   # F: 开始时间

    # Skip blank rows or rows with malformed/non-numeric data (e.g. a
    # corrupted date that can't be parsed as yyyyMMdd) - leave as-is.
    if ($remaining -eq $null -or $totalDays -eq $null -or $startDate -eq $null) { continue }
    if (-not ($startDate -match '^\d{8}$')) { continue }

    $newRemaining = $remaining - 1

    if ($newRemaining -le 0) {
        # Expired today -> renew: restart the clock from today.
        $ws.Cells.Item($r, 6).Value = $today
        $ws.Cells.Item($r, 5).Value = $totalDays
    } else {
        $ws.Cells.Item($r, 5).Value = $newRemaining
    }
}
